$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ A = 'accounts.google.com'; B = "142.251.130.13`n" },
    @{ A = 'www.126.com'; B = "220.181.72.180`n" },
    @{ A = 'urswebzj.nosdn.127.net'; B = "119.188.91.244`n119.188.91.241`n119.167.137.106`n119.167.137.102`n119.167.137.105`n119.188.91.242`n119.188.91.238`n119.188.91.243`n" },
    @{ A = 'mimg.127.net'; B = "61.170.81.231`n61.170.81.215`n61.170.81.250`n" },
    @{ A = 'dl.reg.163.com'; B = "223.252.215.2`n" },
    @{ A = 'passport.126.com'; B = "223.252.215.4`n" },
    @{ A = 'onegoods.nosdn.127.net'; B = "119.188.91.238`n119.188.91.223`n119.188.91.237`n119.188.91.241`n119.167.137.74`n119.167.137.101`n119.188.91.224`n119.167.137.75`n" },
    @{ A = 'utility.mail.163.com'; B = "220.181.12.191`n" },
    @{ A = 'mail.126.com'; B = "123.126.96.204`n" },
    @{ A = 'countly.mail.163.com'; B = "111.124.200.205`n" },
    @{ A = 'cstaticdun.126.net'; B = "183.2.193.244`n183.2.193.238`n14.119.65.239`n" },
    @{ A = 'fl.reg.163.com'; B = "59.111.160.204`n" },
    @{ A = 'content-autofill.googleapis.com'; B = "172.217.27.10`n172.217.27.42`n172.217.31.10`n142.250.199.74`n142.250.204.42`n142.250.204.74`n142.251.130.10`n142.251.222.202`n142.250.66.42`n142.250.66.74`n142.250.66.106`n142.250.66.138`n142.250.207.74`n172.217.24.106`n172.217.24.234`n172.217.25.10`n" },
    @{ A = 'pr.nss.netease.com'; B = "59.111.160.244`n" },
    @{ A = 'mail-activity.nosdn.127.net'; B = "183.2.193.243`n183.2.193.248`n183.2.193.238`n" },
    @{ A = 'b.mail.126.com'; B = "220.181.12.191`n" },
    @{ A = 'mail.163.com'; B = "220.181.12.133`n" },
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item.A
    $ws.Cells.Item($row, 2).Value = $item.B
    $row++
}
